$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.029209278203259
$ws.Cells.Item(2, 4).Value = 1.032308037325799
$ws.Cells.Item(2, 5).Value = 1.042356610010721
$ws.Cells.Item(2, 6).Value = 1.049122976563314
$ws.Cells.Item(2, 9).Value = 1.030352992924349
$ws.Cells.Item(2, 10).Value = 1.034357321054633
$ws.Cells.Item(2, 11).Value = 1.035113643246921
$ws.Cells.Item(2, 12).Value = 1.04513352776347
$ws.Cells.Item(2, 13).Value = 1.051880912504861
$ws.Cells.Item(2, 14).Value = 1.015424632022812
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.030202488257701
$ws.Cells.Item(3, 4).Value = 1.03302903784714
$ws.Cells.Item(3, 5).Value = 1.043371509406735
$ws.Cells.Item(3, 6).Value = 1.050281660440898
$ws.Cells.Item(3, 9).Value = 1.030485586480603
$ws.Cells.Item(3, 10).Value = 1.034991242670671
$ws.Cells.Item(3, 11).Value = 1.035643623758048
$ws.Cells.Item(3, 12).Value = 1.045958707350577
$ws.Cells.Item(3, 13).Value = 1.052850884038371
$ws.Cells.Item(3, 14).Value = 1.015637786353417
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.03084545161165
$ws.Cells.Item(4, 4).Value = 1.033495504600651
$ws.Cells.Item(4, 5).Value = 1.044028928376299
$ws.Cells.Item(4, 6).Value = 1.051032415761432
$ws.Cells.Item(4, 9).Value = 1.030569972184456
$ws.Cells.Item(4, 10).Value = 1.035401140402664
$ws.Cells.Item(4, 11).Value = 1.035985819453958
$ws.Cells.Item(4, 12).Value = 1.046492756877477
$ws.Cells.Item(4, 13).Value = 1.053478950129875
$ws.Cells.Item(4, 14).Value = 1.015775529572001
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.031115822494565
$ws.Cells.Item(5, 4).Value = 1.033691589698016
$ws.Cells.Item(5, 5).Value = 1.044305476978224
$ws.Cells.Item(5, 6).Value = 1.051348274566705
$ws.Cells.Item(5, 9).Value = 1.03060510967482
$ws.Cells.Item(5, 10).Value = 1.035573390981896
$ws.Cells.Item(5, 11).Value = 1.036129501328164
$ws.Cells.Item(5, 12).Value = 1.046717295469083
$ws.Cells.Item(5, 13).Value = 1.053743091532767
$ws.Cells.Item(5, 14).Value = 1.015833393026995
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.031161222996946
$ws.Cells.Item(6, 4).Value = 1.033724512190227
$ws.Cells.Item(6, 5).Value = 1.044351920637505
$ws.Cells.Item(6, 6).Value = 1.051401322787201
$ws.Cells.Item(6, 9).Value = 1.030610989566247
$ws.Cells.Item(6, 10).Value = 1.035602308466655
$ws.Cells.Item(6, 11).Value = 1.036153615735111
$ws.Cells.Item(6, 12).Value = 1.046754997888789
$ws.Cells.Item(6, 13).Value = 1.053787448029703
$ws.Cells.Item(6, 14).Value = 1.015843105980322
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.030849064049572
$ws.Cells.Item(7, 4).Value = 1.033498124769516
$ws.Cells.Item(7, 5).Value = 1.044032622966318
$ws.Cells.Item(7, 6).Value = 1.051036635334691
$ws.Cells.Item(7, 9).Value = 1.030570443022973
$ws.Cells.Item(7, 10).Value = 1.035403442300737
$ws.Cells.Item(7, 11).Value = 1.035987740035369
$ws.Cells.Item(7, 12).Value = 1.04649575707881
$ws.Cells.Item(7, 13).Value = 1.053482479198343
$ws.Cells.Item(7, 14).Value = 1.015776302918714
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.029544877928804
$ws.Cells.Item(8, 4).Value = 1.032551716438082
$ws.Cells.Item(8, 5).Value = 1.042699452397406
$ws.Cells.Item(8, 6).Value = 1.049514349979377
$ws.Cells.Item(8, 9).Value = 1.030398095483686
$ws.Cells.Item(8, 10).Value = 1.034571617959259
$ws.Cells.Item(8, 11).Value = 1.035292904834694
$ws.Cells.Item(8, 12).Value = 1.045412379393043
$ws.Cells.Item(8, 13).Value = 1.052208629698419
$ws.Cells.Item(8, 14).Value = 1.015496706037365
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.027248971704643
$ws.Cells.Item(9, 4).Value = 1.030883543548369
$ws.Cells.Item(9, 5).Value = 1.040355708469733
$ws.Cells.Item(9, 6).Value = 1.046839635700699
$ws.Cells.Item(9, 9).Value = 1.030083607387153
$ws.Cells.Item(9, 10).Value = 1.03310362472123
$ws.Cells.Item(9, 11).Value = 1.034062902700263
$ws.Cells.Item(9, 12).Value = 1.043504142529505
$ws.Cells.Item(9, 13).Value = 1.049967253971992
$ws.Cells.Item(9, 14).Value = 1.015002638541408
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.025719889096323
$ws.Cells.Item(10, 4).Value = 1.029771166583042
$ws.Cells.Item(10, 5).Value = 1.038796925302337
$ws.Cells.Item(10, 6).Value = 1.045061721176115
$ws.Cells.Item(10, 9).Value = 1.029866718497014
$ws.Cells.Item(10, 10).Value = 1.032123505234657
$ws.Cells.Item(10, 11).Value = 1.033239168668101
$ws.Cells.Item(10, 12).Value = 1.042232560466174
$ws.Cells.Item(10, 13).Value = 1.048475252674943
$ws.Cells.Item(10, 14).Value = 1.014672345300848
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.025058143355024
$ws.Cells.Item(11, 4).Value = 1.02928944554938
$ws.Cells.Item(11, 5).Value = 1.0381228419226
$ws.Cells.Item(11, 6).Value = 1.044293107990037
$ws.Cells.Item(11, 9).Value = 1.029771093247723
$ws.Cells.Item(11, 10).Value = 1.031698762919186
$ws.Cells.Item(11, 11).Value = 1.03288160583057
$ws.Cells.Item(11, 12).Value = 1.041682093548769
$ws.Cells.Item(11, 13).Value = 1.047829737449639
$ws.Cells.Item(11, 14).Value = 1.014529110540685
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.024812395324517
$ws.Cells.Item(12, 4).Value = 1.029110505688615
$ws.Cells.Item(12, 5).Value = 1.037872589987436
$ws.Cells.Item(12, 6).Value = 1.044007796632134
$ws.Cells.Item(12, 9).Value = 1.02973531708455
$ws.Cells.Item(12, 10).Value = 1.031540943426437
$ws.Cells.Item(12, 11).Value = 1.032748659399949
$ws.Cells.Item(12, 12).Value = 1.041477646510307
$ws.Cells.Item(12, 13).Value = 1.047590044422325
$ws.Cells.Item(12, 14).Value = 1.014475874590485
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.024865106674181
$ws.Cells.Item(13, 4).Value = 1.029148889207468
$ws.Cells.Item(13, 5).Value = 1.037926263867558
$ws.Cells.Item(13, 6).Value = 1.044068988476106
$ws.Cells.Item(13, 9).Value = 1.029743002802882
$ws.Cells.Item(13, 10).Value = 1.031574798557472
$ws.Cells.Item(13, 11).Value = 1.032777182809066
$ws.Cells.Item(13, 12).Value = 1.041521500142961
$ws.Cells.Item(13, 13).Value = 1.047641455762764
$ws.Cells.Item(13, 14).Value = 1.014487295335487
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.025037828647737
$ws.Cells.Item(14, 4).Value = 1.029274654466398
$ws.Cells.Item(14, 5).Value = 1.038102153311129
$ws.Cells.Item(14, 6).Value = 1.044269520272849
$ws.Cells.Item(14, 9).Value = 1.029768141213186
$ws.Cells.Item(14, 10).Value = 1.031685718554488
$ws.Cells.Item(14, 11).Value = 1.032870619119615
$ws.Cells.Item(14, 12).Value = 1.04166519346458
$ws.Cells.Item(14, 13).Value = 1.047809922714301
$ws.Cells.Item(14, 14).Value = 1.01452471069574
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.025144255490019
$ws.Cells.Item(15, 4).Value = 1.029352141665217
$ws.Cells.Item(15, 5).Value = 1.038210542192227
$ws.Cells.Item(15, 6).Value = 1.044393099157128
$ws.Cells.Item(15, 9).Value = 1.029783595814585
$ws.Cells.Item(15, 10).Value = 1.031754053242086
$ws.Cells.Item(15, 11).Value = 1.032928170885102
$ws.Cells.Item(15, 12).Value = 1.041753730441143
$ws.Cells.Item(15, 13).Value = 1.047913731381152
$ws.Cells.Item(15, 14).Value = 1.014547759275803
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.025763814518463
$ws.Cells.Item(16, 4).Value = 1.029803135756209
$ws.Cells.Item(16, 5).Value = 1.038841680574068
$ws.Cells.Item(16, 6).Value = 1.045112757548197
$ws.Cells.Item(16, 9).Value = 1.029873028813278
$ws.Cells.Item(16, 10).Value = 1.032151686768545
$ws.Cells.Item(16, 11).Value = 1.033262880436679
$ws.Cells.Item(16, 12).Value = 1.042269096050613
$ws.Cells.Item(16, 13).Value = 1.04851810460841
$ws.Cells.Item(16, 14).Value = 1.014681846795461
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.026152543243296
$ws.Cells.Item(17, 4).Value = 1.030086018532436
$ws.Cells.Item(17, 5).Value = 1.039237812917293
$ws.Cells.Item(17, 6).Value = 1.045564511155233
$ws.Cells.Item(17, 9).Value = 1.029928669960309
$ws.Cells.Item(17, 10).Value = 1.032401019937994
$ws.Cells.Item(17, 11).Value = 1.033472599582852
$ws.Cells.Item(17, 12).Value = 1.04259240809302
$ws.Cells.Item(17, 13).Value = 1.048897354365065
$ws.Cells.Item(17, 14).Value = 1.014765898777983
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.026379316578273
$ws.Cells.Item(18, 4).Value = 1.030251013950069
$ws.Cells.Item(18, 5).Value = 1.039468955090899
$ws.Cells.Item(18, 6).Value = 1.045828130590333
$ws.Cells.Item(18, 9).Value = 1.029960959396042
$ws.Cells.Item(18, 10).Value = 1.032546418527935
$ws.Cells.Item(18, 11).Value = 1.033594840209065
$ws.Cells.Item(18, 12).Value = 1.042781003683006
$ws.Cells.Item(18, 13).Value = 1.049118615565139
$ws.Cells.Item(18, 14).Value = 1.014814904060294
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.026456646263408
$ws.Cells.Item(19, 4).Value = 1.030307272191975
$ws.Cells.Item(19, 5).Value = 1.039547782982808
$ws.Cells.Item(19, 6).Value = 1.045918038218102
$ws.Cells.Item(19, 9).Value = 1.029971941243759
$ws.Cells.Item(19, 10).Value = 1.032595990040821
$ws.Cells.Item(19, 11).Value = 1.033636506666591
$ws.Cells.Item(19, 12).Value = 1.042845312123784
$ws.Cells.Item(19, 13).Value = 1.049194068598533
$ws.Cells.Item(19, 14).Value = 1.014831610057833
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.026110832765877
$ws.Cells.Item(20, 4).Value = 1.030055668430963
$ws.Cells.Item(20, 5).Value = 1.039195302883097
$ws.Cells.Item(20, 6).Value = 1.045516029937846
$ws.Cells.Item(20, 9).Value = 1.029922717266673
$ws.Cells.Item(20, 10).Value = 1.032374272282607
$ws.Cells.Item(20, 11).Value = 1.033450107491
$ws.Cells.Item(20, 12).Value = 1.042557718408162
$ws.Cells.Item(20, 13).Value = 1.048856659154048
$ws.Cells.Item(20, 14).Value = 1.014756882947389
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.024986964845089
$ws.Cells.Item(21, 4).Value = 1.02923761992385
$ws.Cells.Item(21, 5).Value = 1.038050354588768
$ws.Cells.Item(21, 6).Value = 1.044210463526348
$ws.Cells.Item(21, 9).Value = 1.02976074566086
$ws.Cells.Item(21, 10).Value = 1.031653056792369
$ws.Cells.Item(21, 11).Value = 1.032843108084797
$ws.Cells.Item(21, 12).Value = 1.041622878783291
$ws.Cells.Item(21, 13).Value = 1.047760311183563
$ws.Cells.Item(21, 14).Value = 1.014513693691416
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.024280655836066
$ws.Cells.Item(22, 4).Value = 1.02872323906031
$ws.Cells.Item(22, 5).Value = 1.037331247816246
$ws.Cells.Item(22, 6).Value = 1.043390677189009
$ws.Cells.Item(22, 9).Value = 1.029657422745232
$ws.Cells.Item(22, 10).Value = 1.031199303474358
$ws.Cells.Item(22, 11).Value = 1.032460701959304
$ws.Cells.Item(22, 12).Value = 1.041035229024694
$ws.Cells.Item(22, 13).Value = 1.04707145708849
$ws.Cells.Item(22, 14).Value = 1.014360604858359
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.024655054041482
$ws.Cells.Item(23, 4).Value = 1.028995925682192
$ws.Cells.Item(23, 5).Value = 1.037712387011531
$ws.Cells.Item(23, 6).Value = 1.043825159471691
$ws.Cells.Item(23, 9).Value = 1.029712336812496
$ws.Cells.Item(23, 10).Value = 1.031439874661302
$ws.Cells.Item(23, 11).Value = 1.032663494656014
$ws.Cells.Item(23, 12).Value = 1.041346741683404
$ws.Cells.Item(23, 13).Value = 1.047436587664689
$ws.Cells.Item(23, 14).Value = 1.014441777715
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.02612967982868
$ws.Cells.Item(24, 4).Value = 1.030069382351906
$ws.Cells.Item(24, 5).Value = 1.039214511075412
$ws.Cells.Item(24, 6).Value = 1.045537936144721
$ws.Cells.Item(24, 9).Value = 1.029925407542855
$ws.Cells.Item(24, 10).Value = 1.032386358499384
$ws.Cells.Item(24, 11).Value = 1.033460270961804
$ws.Cells.Item(24, 12).Value = 1.042573393142882
$ws.Cells.Item(24, 13).Value = 1.048875047409818
$ws.Cells.Item(24, 14).Value = 1.014760956877476
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.027842251640741
$ws.Cells.Item(25, 4).Value = 1.031314856836812
$ws.Cells.Item(25, 5).Value = 1.040960970480536
$ws.Cells.Item(25, 6).Value = 1.047530193380414
$ws.Cells.Item(25, 9).Value = 1.030166186238466
$ws.Cells.Item(25, 10).Value = 1.033483395239448
$ws.Cells.Item(25, 11).Value = 1.034381548687026
$ws.Cells.Item(25, 12).Value = 1.043997368147063
$ws.Cells.Item(25, 13).Value = 1.050546308643906
$ws.Cells.Item(25, 14).Value = 1.015130529137956
